# Fill in diary entries for rows 26 and 27 (weeks of 3/5/2020 and 3/12/2020)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Match styling of the preceding filled rows (e.g. row 25) for the two new rows first,
# so the new text inherits the same wrap/format as the rest of the diary table.
$ws.Range("A25:G25").Copy()
$ws.Range("A26:G27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Dates and the reused Time/Participants/Mood values for both rows
$ws.Range("A26").Value = [DateTime]"2020-03-05"
$ws.Range("B26").Value = "5:00 -7:50 pm"
$ws.Range("C26").Value = "N/A"
$ws.Range("G26").Value = "Good"

$ws.Range("A27").Value = [DateTime]"2020-03-12"
$ws.Range("B27").Value = "5:00 -7:50 pm"
$ws.Range("C27").Value = "N/A"
$ws.Range("G27").Value = "Good"

# Goal column for both new rows (entered first)
$ws.Range("D26").Value = "Testing"
$ws.Range("D27").Value = "Advancd topics"

# Achievements + Reflection for row 26
$ws.Range("E26").Value = "Learned how testing can be used to gain an undestanding of codebase"
$ws.Range("F26").Value = "We learned a ton of concepts in Jones' class about testing and now that I'm familiar with what to look for….actaully doing the Pacman activitiy and reading the test cases really made me realize the value testing has on understanding key componenets to software. It was great because this is probably the first thing Im gonig to do in my internship to learn about the compan'ys codebase. Also, the guest speakers were awesome. They made me want to really explore my passion and I hope one day to be as inspirational as them"

# Achievements + Reflection for row 27
$ws.Range("E27").Value = "Important to keep a good attitude and to always learn to stay releavant and advance in your career"
$ws.Range("F27").Value = "Jve learned a lot from this class and a key takeway for me is to always keep learning new things, perspectives, mental models. Im inspired to continue to learn new ways to challenge myself. Ill probably reach out to Andre or Kaj in the future about what are good resources to keep learning about software"

# Row heights grow to fit the new wrapped text
$ws.Rows.Item(26).RowHeight = 238
$ws.Rows.Item(27).RowHeight = 153

# Update view: scroll position and active selection to match the edited author's final view
$excel.ActiveWindow.ScrollRow = 24
$ws.Range("B26").Select()
